$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = 45981
$ws.Range("A32").NumberFormat = $ws.Range("A31").NumberFormat
$ws.Range("B32").Value = 68
$ws.Range("C32").Value = 80
$ws.Range("D32").Value = 78
